# Applies the 2024-08-18 cryptos-list refresh: updated prices/volume
# percentages, plus the Binance-PegBSC-USD / Kaspa row-order swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.751.25"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.619.09"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'532.43"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "'142.90"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.567"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'6.90"
$ws.Range("E9").Value = "  +7.12%  "
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").Value = "'0.335"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D13").Value = "3.086.30"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "58.679.17"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "'20.88"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "2.618.78"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "'4.40"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "'335.57"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").Value = "'10.14"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'6.21"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'66.32"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("D24").Value = "'0.413"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.163"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("D27").Value = "'7.18"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").Value = "0.0₃0734"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'1.63"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").Value = "'5.81"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'151.11"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "'18.67"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").Value = "'3.93"
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "'0.829"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "'0.821"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'281.68"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "'0.593"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").Value = "'10.72"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'0.0533"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D46").Value = "'18.87"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").Value = "'0.0224"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "1.940.39"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "'17.97"
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").Value = "'111.62"
$ws.Range("E51").Value = "  +0.63%  "
